# Insert a new record (row) for Jengibre / Vega Modelo de Temuco just before
# the existing row 139, shifting all subsequent rows down by one, and fill in
# the new row with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 139 (and everything below it) down by one row.
$ws.Rows(139).Insert()

# Populate the newly inserted row 139 with the new weekly data point.
$ws.Cells.Item(139, 1).Value = 10
$ws.Cells.Item(139, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(139, 3).Value = "La Araucanía"
$ws.Cells.Item(139, 4).Value = 44879
$ws.Cells.Item(139, 5).Value = 9
$ws.Cells.Item(139, 6).Value = 100114007
$ws.Cells.Item(139, 7).Value = "Jengibre"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 90
$ws.Cells.Item(139, 11).Value = 18000
$ws.Cells.Item(139, 12).Value = 20000
$ws.Cells.Item(139, 13).Value = 18889
$ws.Cells.Item(139, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(139, 15).Value = "Perú"
$ws.Cells.Item(139, 16).Value = 1453
$ws.Cells.Item(139, 17).Value = 13
$ws.Cells.Item(139, 18).Value = "Hortaliza"
